$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for the columns that get permuted
# across rows: D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
$rows = @(2,3,4,5,6,7,10,11,12,13,14,15,16,17,18,19,20,21)

$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
}

# Mapping: target row -> source row (whose pre-edit values get written into target)
$mapping = @{
    2  = 20
    3  = 21
    4  = 13
    5  = 14
    6  = 2
    7  = 3
    10 = 18
    11 = 16
    12 = 17
    13 = 11
    14 = 12
    15 = 7
    16 = 19
    17 = 10
    18 = 6
    19 = 4
    20 = 5
    21 = 15
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $vals = $snapshot[$source]

    $ws.Cells.Item($target, 4).Value2 = $vals.D
    $ws.Cells.Item($target, 12).Value2 = $vals.L
    $ws.Cells.Item($target, 13).Value2 = $vals.M
    $ws.Cells.Item($target, 14).Value2 = $vals.N
    $ws.Cells.Item($target, 15).Value2 = $vals.O
    $ws.Cells.Item($target, 16).Value2 = $vals.P
    $ws.Cells.Item($target, 19).Value2 = $vals.S
}
